$wb = $excel.ActiveWorkbook

# Germany sheet is the template used for the new Italy sheet.
$germany = $wb.Worksheets.Item("Germany")
$slovakia = $wb.Worksheets.Item("Slovakia")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy Germany to the end of the workbook and rename it to Italy.
$germany.Copy($null, $lastSheet)
$italy = $wb.Worksheets.Item($wb.Worksheets.Count)
$italy.Name = "Italy"

# Fill in the Italy-specific market name and part numbers.
$italy.Range("B2").Value = "Italy Market"
$italy.Range("B4").Value = "NGC-3145/T2454/T2453/T2452/T2455"

# Slovakia is no longer the tab that was last interacted with - move its
# selection without leaving it as the active sheet.
$slovakia.Activate()
$slovakia.Range("F10").Select()

# Germany's selection becomes a full-column selection.
$germany.Activate()
$germany.Range("A1:XFD1048576").Select()

# Italy ends up as the active sheet with B4 selected.
$italy.Activate()
$italy.Range("B4").Select()
